# Apply PyST development plan updates: support observation groups and
# parameter groups in the PEST control file object.

$wb = $excel.ActiveWorkbook

$wsParserWriter = $wb.Worksheets.Item("ParserWriter")
$wsPstFile = $wb.Worksheets.Item("PST file")

# --- Update "PST file" sheet data ---------------------------------------
# block separation -> complete
$wsPstFile.Range("B5").Value = 1
# control data -> complete
$wsPstFile.Range("B7").Value = 1
# parameter groups -> complete, implemented via BeoJACTEST
$wsPstFile.Range("B13").Value = 1
$wsPstFile.Range("C13").Value = "BeoJACTEST"
# parameter data -> complete, implemented via BeoJACTEST
$wsPstFile.Range("B14").Value = 1
$wsPstFile.Range("C14").Value = "BeoJACTEST"
# observation groups -> complete
$wsPstFile.Range("B15").Value = 1

# --- Update sheet selections ---------------------------------------------
$wsPstFile.Range("B16").Select()
$wsParserWriter.Range("C6").Select()

# --- Make ParserWriter the active sheet/tab ------------------------------
$wsParserWriter.Activate()

# Recalculate all formulas so the AVERAGE() on ParserWriter reflects the
# updated PST file completion values.
$excel.Calculate()
